$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) format on the Price cells whose new values would
# otherwise be auto-coerced into numbers by Excel (e.g. "593.41", "1.00"),
# matching the original inline-string (text) storage of column D.
foreach ($addr in @('D5', 'D6', 'D7', 'D9', 'D11', 'D12', 'D13', 'D14', 'D16', 'D19', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D31', 'D33', 'D34', 'D35', 'D36', 'D46', 'D47', 'D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptos list values (price + 1h volume change),
# plus the ThetaToken/dogwifhat row swap (rows 43-44).
$ws.Range('D2').Value = '70.725.61'
$ws.Range('E2').Value = '  +5.73%  '
$ws.Range('D3').Value = '3.636.10'
$ws.Range('E3').Value = '  +5.69%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '593.41'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').Value = '195.14'
$ws.Range('E6').Value = '  +3.37%  '
$ws.Range('D7').Value = '0.645'
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('D8').Value = '3.628.40'
$ws.Range('E8').Value = '  +5.71%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D11').Value = '0.680'
$ws.Range('E11').Value = '  +5.43%  '
$ws.Range('D12').Value = '58.11'
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').Value = '0.0000315'
$ws.Range('E13').Value = '  +14.16%  '
$ws.Range('D14').Value = '9.98'
$ws.Range('E14').Value = '  +5.50%  '
$ws.Range('D15').Value = '4.220.77'
$ws.Range('E15').Value = '  +5.96%  '
$ws.Range('D16').Value = '20.58'
$ws.Range('E16').Value = '  +9.29%  '
$ws.Range('D17').Value = '3.640.02'
$ws.Range('E17').Value = '  +5.83%  '
$ws.Range('D18').Value = '70.771.40'
$ws.Range('E18').Value = '  +5.87%  '
$ws.Range('D19').Value = '12.78'
$ws.Range('E19').Value = '  +5.90%  '
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('E21').Value = '  +4.22%  '
$ws.Range('D22').Value = '489.26'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = '19.33'
$ws.Range('E23').Value = '  +13.22%  '
$ws.Range('D24').Value = '5.26'
$ws.Range('E24').Value = '  -2.99%  '
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('D26').Value = '91.44'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').Value = '3.18'
$ws.Range('E27').Value = '  +6.52%  '
$ws.Range('D28').Value = '11.46'
$ws.Range('E28').Value = '  +4.77%  '
$ws.Range('E29').Value = '  +6.79%  '
$ws.Range('E30').Value = '  +7.11%  '
$ws.Range('D31').Value = '32.84'
$ws.Range('E31').Value = '  +5.64%  '
$ws.Range('E32').Value = '  +10.04%  '
$ws.Range('D33').Value = '12.30'
$ws.Range('E33').Value = '  +4.65%  '
$ws.Range('D34').Value = '66.34'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('D35').Value = '612.35'
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('D36').Value = '40.63'
$ws.Range('E36').Value = '  +9.68%  '
$ws.Range('D37').Value = '0.0₃0838'
$ws.Range('E37').Value = '  +11.61%  '
$ws.Range('E38').Value = '  +6.04%  '
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E41').Value = '  +2.77%  '
$ws.Range('D42').Value = '3.326.21'
$ws.Range('E42').Value = '  +3.90%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E43').Value = '  +17.13%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('E44').Value = '  +9.64%  '
$ws.Range('E45').Value = '  +9.77%  '
$ws.Range('D46').Value = '0.0461'
$ws.Range('E46').Value = '  +7.20%  '
$ws.Range('D47').Value = '9.70'
$ws.Range('E47').Value = '  +13.08%  '
$ws.Range('E48').Value = '  +2.99%  '
$ws.Range('E49').Value = '  +3.65%  '
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.11%  '
